# Update of the TODO list - add tasks 89 and 90 to the "short term" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sanity: make sure we are editing the "short term" sheet (it is the active
# sheet/tab in the workbook already, but select it explicitly to be safe).
$ws = $wb.Worksheets.Item("short term")

$orange = 49407   # RGB(255,192,0) == 0xFFC000, same amber used elsewhere in the sheet

# --- Row 22: task 89 -------------------------------------------------
$a22 = $ws.Range("A22")
$a22.Value = "89. change the default for log scale in each plot (fit or BMD) - log if max/min > 100 and get this default choice in the Shiny app."
$a22.Interior.Color = $orange
$a22.Borders.LineStyle = 1
$a22.WrapText = $true

$b22 = $ws.Range("B22")
$b22.Value = "ML and A"
$b22.Interior.Color = $orange

# --- Row 23: task 90 -------------------------------------------------
$a23 = $ws.Range("A23")
$a23.Value = "90. add an explaination of minBMD and or a option to define it in Shiny"
$a23.Interior.Color = $orange
$a23.Borders.LineStyle = 1
$a23.WrapText = $true

$b23 = $ws.Range("B23")
$b23.Value = "ML and A"
$b23.Interior.Color = $orange

# Match the saved selection state recorded in the workbook after the edit.
$ws.Range("A27").Select() | Out-Null
